# Apply weekly price update: insert a new data row for Acelga (Femacal de La Calera)
# at worksheet row 530, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 530 (pushes old row 530..630 down to 531..631)
$ws.Range("A530:R530").EntireRow.Insert()

# The row that used to be 530 is now at 531. Copy its "static" (non price/date)
# column values into the freshly inserted row 530, since every record in this
# sheet shares the same Mercado/Region/Categoria/Unidad/Origen/etc.
$staticCols = @("A","B","C","E","F","G","H","I","N","O","Q","R")
foreach ($col in $staticCols) {
    $ws.Range($col + "530").Value = $ws.Range($col + "531").Value2
}

# Match the date number format used by the other rows in column D
$ws.Range("D530").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new record's date/volume/price values
$ws.Range("D530").Value = 45209
$ws.Range("J530").Value = 230
$ws.Range("K530").Value = 3000
$ws.Range("L530").Value = 3500
$ws.Range("M530").Value = 3261
$ws.Range("P530").Value = 544
